# Insert a new data row at row 491 (pushing existing rows 491:527 down to
# 492:528) and populate it with the new weekly price entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(491).Insert()

$ws.Range("A491").Value = 5
$ws.Range("B491").Value = "Macroferia Regional de Talca"
$ws.Range("C491").Value = "Maule"
$ws.Range("D491").Value = 45013
$ws.Range("E491").Value = 7
$ws.Range("F491").Value = 100112032
$ws.Range("G491").Value = "Zapallo italiano"
$ws.Range("H491").Value = "Sin especificar"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 300
$ws.Range("K491").Value = 5000
$ws.Range("L491").Value = 5000
$ws.Range("M491").Value = 5000
$ws.Range("N491").Value = "$/caja 50 unidades"
$ws.Range("O491").Value = "Región del Maule"
$ws.Range("P491").Value = 100
$ws.Range("Q491").Value = 50
$ws.Range("R491").Value = "Hortaliza"
